# Add a new "testRow" entry as row 18 of the TestData sheet
# (new unit-test fixture row), then move the selection onto it.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")
$ws.Activate()

$ws.Range("A18").Value = "testRowValue"
$ws.Range("B18").Value = "testRowData."

$ws.Range("B18").Select()
